# Springboot-projects-repository.xlsx -- "Add files via upload"
# Adds two new project rows (React+Spring Boot CRUD / React+Spring Boot
# CRUD-without-DB) to the tracking sheet, numbering the existing last row
# and filling in the two new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# 1. Row 18 gains its running S.No (16) in column A.
# ---------------------------------------------------------------------
$ws.Range("A18").Value() = 16

# ---------------------------------------------------------------------
# 2. D15 / D16 switch to the (equivalent) centered+wrapped border style
#    that D11/D18 already use.
# ---------------------------------------------------------------------
Copy-CellFormat "D11" "D15"
Copy-CellFormat "D11" "D16"

# ---------------------------------------------------------------------
# 3. Row 19: was a blank placeholder row -- turn it into project #17
#    (Reactjs + Spring Boot CRUD Full Stack Application).
# ---------------------------------------------------------------------
Copy-CellFormat "A2"  "B19"
Copy-CellFormat "A2"  "C19"
Copy-CellFormat "O17" "O19"

$ws.Range("A19").Value() = 17
$ws.Range("B19").Value() = "Reactjs + Spring Boot CRUD Full Stack Application"
$ws.Range("C19").Value() = "React - JS, Springboot, REST, MySQL"
$ws.Range("D19").Value() = "D:\springboot-development\springbootreactmysql"
$ws.Range("E19").Value() = "SpringToolSuite4"
$ws.Range("F19").Value() = "springbootreactmysql"
$ws.Range("G19").Value() = "Spring Starter Project"
$ws.Range("H19").Value() = 17
$ws.Range("I19").Value() = "Yes"
$ws.Range("J19").Value() = "Yes"
$ws.Range("K19").Value() = "Spring Web"
$ws.Range("L19").Value() = "Spring Boot DevTools"
$ws.Range("M19").Value() = "Spring Data JPA"
$ws.Range("N19").Value() = "MySQL"
$ws.Range("O19").Value() = "Lombok and Spring Security"
$ws.Range("Q19").Value() = "yes"
$ws.Range("R19").Value() = "springboot-react-mysql-CURD"

$ws.Rows.Item(19).RowHeight = 15

# ---------------------------------------------------------------------
# 4. Row 20: brand-new project #18
#    (Reactjs + Spring Boot -User CURD integration without DB).
# ---------------------------------------------------------------------
Copy-CellFormat "A2"  "A20"
Copy-CellFormat "A2"  "B20"
Copy-CellFormat "A2"  "C20"
Copy-CellFormat "B2"  "D20"
Copy-CellFormat "A2"  "E20"
Copy-CellFormat "B2"  "F20"
Copy-CellFormat "A2"  "G20"
Copy-CellFormat "A2"  "H20"
Copy-CellFormat "A2"  "I20"
Copy-CellFormat "A2"  "J20"
Copy-CellFormat "A2"  "K20"
Copy-CellFormat "A2"  "L20"
Copy-CellFormat "A2"  "M20"
Copy-CellFormat "A2"  "N20"
Copy-CellFormat "O17" "O20"
Copy-CellFormat "A2"  "P20"
Copy-CellFormat "A2"  "Q20"
Copy-CellFormat "A2"  "R20"
Copy-CellFormat "B2"  "S20"

$ws.Range("A20").Value() = 18
$ws.Range("B20").Value() = "Reactjs + Spring Boot -User CURD integration without DB"
$ws.Range("C20").Value() = "React - JS, Springboot, REST"
$ws.Range("D20").Value() = "// React code - D:\react-development\react-springboot-sathees\my-test-get`n// Springboot code - D:\springboot-development\reactspringbootwebsevice"
$ws.Range("E20").Value() = "SpringToolSuite4"
$ws.Range("F20").Value() = "reactspringbootwebsevice"
$ws.Range("G20").Value() = "Spring Starter Project"
$ws.Range("H20").Value() = 17
$ws.Range("I20").Value() = "Yes"
$ws.Range("J20").Value() = "Yes"
$ws.Range("K20").Value() = "Spring Web"
$ws.Range("L20").Value() = "Spring Boot DevTools"
$ws.Range("M20").Value() = "Spring Data JPA"
$ws.Range("N20").Value() = "Postgres SQL Driver"
$ws.Range("O20").Value() = "Lombok"
$ws.Range("R20").Value() = "springboot-react-User-CURD"

$ws.Rows.Item(20).RowHeight = 58.2

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5. Selection ends on the new last cell, same as the saved workbook.
# ---------------------------------------------------------------------
$ws.Range("S20").Select()
